$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) to the table
$ws.Range("A5").Value = "Saccharomyceo cochise"
$ws.Range("B5").Value = "S_cerevisiae_16S.fasta"
$ws.Range("C5").Value = "B_subtilis_recA.fasta"
$ws.Range("D5").Value = "S_cerevisiae_rpoB.fasta"

# Set print page orientation to portrait (adds pageSetup element)
$ws.PageSetup.Orientation = 1

# Move the active selection to E12, mirroring the saved cursor position
$ws.Range("E12").Select()

$wb.Save()
